# Halotype 1a: mark plate 011 as DONE and record the "plate needed" tally
# row for plate 012 (the newly finished plate), then move the active
# selection onto plate 012's sheet.

$wb = $excel.ActiveWorkbook

# 1) "PocHistone RLFP 011" is finished -> prefix its tab name with "DONE ".
$ws011 = $wb.Worksheets.Item("PocHistone RLFP 011")
$ws011.Name = "DONE PocHistone RLFP 011"

# 2) Add the "plate needed" tally row (row 11) to "PocHistone RLFP 012",
#    matching the same row shape used on the other completed plates.
$ws012 = $wb.Worksheets.Item("PocHistone RLFP 012")
$ws012.Range("A11").Value = "plate needed"
$ws012.Range("B11").Value = 18
$ws012.Range("C11").Value = 17
$ws012.Range("D11").Value = 28
$ws012.Range("E11").Value = 29
$ws012.Range("F11").Value = 20
$ws012.Range("G11").Value = 25
$ws012.Range("H11").Value = 30

# 3) Select the new active cell on plate 012 and make it the active sheet.
$ws012.Range("I25").Select()
$ws012.Activate()
